# Quarterly indexing bug-fix: the forecast-error table holds one row of
# metrics (ME, MAE, MSE, RMSE, SE, N) per evaluated quarter, most recent
# quarter on top (row 2). A new quarter's evaluation needs to be inserted
# at the top, so every existing quarter's row of results shifts down by
# one row, and the oldest quarter (previously in row 11 / Q9) rolls off
# the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing metric rows (2-10) down into (3-11). The last row's
# previous contents (old row 11 / Q9) are overwritten and discarded.
$ws.Range("B2:G10").Copy()
$ws.Range("B3:G11").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Write the newly computed metrics for the latest quarter into row 2.
$ws.Range("B2").Value = 0.1724578193461484
$ws.Range("C2").Value = 0.39058239716261
$ws.Range("D2").Value = 0.3033305724894426
$ws.Range("E2").Value = 0.550754548314803
$ws.Range("F2").Value = 0.5414156770869448
$ws.Range("G2").Value = 15
